$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.863.20"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "2.638.30"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.86"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.42"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.57"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "3.106.36"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.25"
$ws.Range("E14").Value = "  +6.06%  "
$ws.Range("D15").Value = "60.828.77"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "2.648.68"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.11"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.83"
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.40"
$ws.Range("E27").Value = "  +5.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.00"
$ws.Range("E28").Value = "  +6.41%  "
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("E30").Value = "  +5.90%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.05"
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.08"
$ws.Range("E34").Value = "  +9.52%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.59"
$ws.Range("E35").Value = "  +8.05%  "
$ws.Range("E36").Value = "  +7.15%  "
$ws.Range("E37").Value = "  +3.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "341.52"
$ws.Range("E38").Value = "  +9.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.12"
$ws.Range("E39").Value = "  +5.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.905"
$ws.Range("E40").Value = "  +7.41%  "
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.31"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.77"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0574"
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.96"
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("E46").Value = "  +3.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.28"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("E48").Value = "  +3.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0994"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "2.086.70"
$ws.Range("E51").Value = "  +2.55%  "
